# Admin tasks + new code for frequency correlation in SISO systems
#
# Highlight the "2 . Imperfect CSI dérivation -> présenter les résultats et
# montrer à Philippe" to-do item in red so it stands out among the other
# priority items.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive text instead of a hard
# coded index, so the edit is resilient to unrelated paragraph insertions
# elsewhere in the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*CSI d*rivation*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Imperfect CSI derivation' paragraph"
}

# Paint the whole paragraph (all runs + the paragraph mark) red
# (wdColorRed = 255 => RGB FF0000), matching Word's "font color" UI action.
$target.Range.Font.Color = 255

Write-Output ("Recolored paragraph: " + $target.Range.Text)
